$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1743.6666
$ws.Range("I43").Value = 2445
$ws.Range("J43").Value = 1543.2858
$ws.Range("K43").Value = 2445
$ws.Range("L43").Value = 1543.2858
$ws.Range("M43").Value = -2376
$ws.Range("N43").Value = -1681.2858
$ws.Range("H107").Value = 5121.864
$ws.Range("I107").Value = 5909.2104
$ws.Range("K107").Value = 5909.2104
$ws.Range("M107").Value = -3989.2104
$ws.Range("H112").Value = 18337.266
$ws.Range("J112").Value = 19790.92
$ws.Range("L112").Value = 59372.75999999999
$ws.Range("N112").Value = -61588.75999999999
$ws.Range("H118").Value = 787.4666999999999
$ws.Range("I118").Value = 613.5833
$ws.Range("K118").Value = 1840.7499
$ws.Range("M118").Value = -183.7499
$ws.Range("H138").Value = 9826.527
$ws.Range("J138").Value = 9169.821
$ws.Range("L138").Value = 27509.463
$ws.Range("N138").Value = -37789.463

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 17500
$ws.Range("J25").Value = 17500
$ws.Range("L25").Value = 17500
$ws.Range("N25").Value = -18304
$ws.Range("H32").Value = 1589.4166
$ws.Range("I32").Value = 1647.2122
$ws.Range("J32").Value = 953.6667
$ws.Range("K32").Value = 1647.2122
$ws.Range("L32").Value = 953.6667
$ws.Range("M32").Value = -1360.2122
$ws.Range("N32").Value = -1527.6667
$ws.Range("H74").Value = 46111
$ws.Range("I74").Value = 52301.895
$ws.Range("J74").Value = 2774.75
$ws.Range("K74").Value = 52301.895
$ws.Range("L74").Value = 2774.75
$ws.Range("M74").Value = -51427.895
$ws.Range("N74").Value = -4522.75
$ws.Range("H77").Value = 46111
$ws.Range("I77").Value = 52301.895
$ws.Range("J77").Value = 2774.75
$ws.Range("K77").Value = 261509.475
$ws.Range("L77").Value = 13873.75
$ws.Range("M77").Value = -257141.475
$ws.Range("N77").Value = -22609.75
$ws.Range("H92").Value = 220029900
$ws.Range("I92").Value = 90000
$ws.Range("J92").Value = 275014880
$ws.Range("K92").Value = 90000
$ws.Range("L92").Value = 275014880
$ws.Range("M92").Value = -87504
$ws.Range("N92").Value = -275019872
$ws.Range("H124").Value = 61500
$ws.Range("J124").Value = 61500
$ws.Range("L124").Value = 61500
$ws.Range("N124").Value = -71320
$ws.Range("H132").Value = 2059.408
$ws.Range("I132").Value = 1517.5938
$ws.Range("K132").Value = 4552.7814
$ws.Range("M132").Value = -2022.7814

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2704.2273
$ws.Range("I20").Value = 1565.0834
$ws.Range("J20").Value = 4071.2
$ws.Range("K20").Value = 1565.0834
$ws.Range("L20").Value = 4071.2
$ws.Range("M20").Value = -1318.0834
$ws.Range("N20").Value = -4565.2
$ws.Range("H32").Value = 52994.668
$ws.Range("J32").Value = 52994.668
$ws.Range("L32").Value = 52994.668
$ws.Range("N32").Value = -53762.668
$ws.Range("H82").Value = 31001.363
$ws.Range("J82").Value = 70919.664
$ws.Range("L82").Value = 70919.664
$ws.Range("N82").Value = -71685.664
$ws.Range("H85").Value = 31001.363
$ws.Range("J85").Value = 70919.664
$ws.Range("L85").Value = 70919.664
$ws.Range("N85").Value = -73571.664
$ws.Range("H107").Value = 3111.2
$ws.Range("I107").Value = 3099.3
$ws.Range("K107").Value = 3099.3
$ws.Range("M107").Value = -1179.3
$ws.Range("H134").Value = 6812.3706
$ws.Range("I134").Value = 7150.8696
$ws.Range("K134").Value = 21452.6088
$ws.Range("M134").Value = -18917.6088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1215.2
$ws.Range("I22").Value = 1026.1818
$ws.Range("J22").Value = 1363.7142
$ws.Range("K22").Value = 1026.1818
$ws.Range("L22").Value = 1363.7142
$ws.Range("M22").Value = -676.1818000000001
$ws.Range("N22").Value = -2063.7142
$ws.Range("H74").Value = 59165
$ws.Range("I74").Value = 58499
$ws.Range("J74").Value = 59498
$ws.Range("K74").Value = 58499
$ws.Range("L74").Value = 59498
$ws.Range("M74").Value = -57625
$ws.Range("N74").Value = -61246
$ws.Range("H77").Value = 59165
$ws.Range("I77").Value = 58499
$ws.Range("J77").Value = 59498
$ws.Range("K77").Value = 175497
$ws.Range("L77").Value = 178494
$ws.Range("M77").Value = -171129
$ws.Range("N77").Value = -187230
$ws.Range("H107").Value = 58836956
$ws.Range("I107").Value = 90929220
$ws.Range("J107").Value = 1145.1666
$ws.Range("K107").Value = 90929220
$ws.Range("L107").Value = 1145.1666
$ws.Range("M107").Value = -90927300
$ws.Range("N107").Value = -4985.1666
$ws.Range("H132").Value = 21041.186
$ws.Range("I132").Value = 7187.091
$ws.Range("K132").Value = 21561.273
$ws.Range("M132").Value = -19031.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 6.6
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 6.6
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = 19.8
$ws.Range("N12").Value = -365.8
$ws.Range("H64").Value = 2859.4
$ws.Range("J64").Value = 2824.5
$ws.Range("L64").Value = 8473.5
$ws.Range("N64").Value = -9013.5
$ws.Range("H67").Value = 2859.4
$ws.Range("J67").Value = 2824.5
$ws.Range("L67").Value = 8473.5
$ws.Range("N67").Value = -10345.5
$ws.Range("H117").Value = 9291.066000000001
$ws.Range("I117").Value = 3712.25
$ws.Range("J117").Value = 11319.728
$ws.Range("K117").Value = 11136.75
$ws.Range("L117").Value = 33959.18399999999
$ws.Range("M117").Value = -7694.75
$ws.Range("N117").Value = -40843.18399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -21166
$ws.Range("H80").Value = 3716.647
$ws.Range("I80").Value = 3286
$ws.Range("K80").Value = 3286
$ws.Range("M80").Value = -2288
$ws.Range("H83").Value = 3716.647
$ws.Range("I83").Value = 3286
$ws.Range("K83").Value = 16430
$ws.Range("M83").Value = -11438
$ws.Range("H107").Value = 725.25
$ws.Range("J107").Value = 737.6667
$ws.Range("L107").Value = 737.6667
$ws.Range("N107").Value = -4577.6667
$ws.Range("H113").Value = 2760.0952
$ws.Range("I113").Value = 2346.9
$ws.Range("K113").Value = 2346.9
$ws.Range("M113").Value = -176.9000000000001
$ws.Range("H122").Value = 7582.325
$ws.Range("I122").Value = 9221.634
$ws.Range("K122").Value = 27664.902
$ws.Range("M122").Value = -25214.902
$ws.Range("H132").Value = 4085.1667
$ws.Range("I132").Value = 3422.5833
$ws.Range("J132").Value = 6735.5
$ws.Range("K132").Value = 10267.7499
$ws.Range("L132").Value = 20206.5
$ws.Range("M132").Value = -7737.749899999999
$ws.Range("N132").Value = -25266.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58392.75
$ws.Range("J7").Value = 9761.75
$ws.Range("L7").Value = 9761.75
$ws.Range("N7").Value = -9985.75
$ws.Range("H46").Value = 1937.5
$ws.Range("I46").Value = 1782.5
$ws.Range("J46").Value = 2040.8334
$ws.Range("K46").Value = 1782.5
$ws.Range("L46").Value = 2040.8334
$ws.Range("M46").Value = -1594.5
$ws.Range("N46").Value = -2416.8334
$ws.Range("H59").Value = 57950
$ws.Range("J59").Value = 57950
$ws.Range("L59").Value = 57950
$ws.Range("N59").Value = -59258
$ws.Range("H122").Value = 4724.0356
$ws.Range("I122").Value = 3662.9
$ws.Range("J122").Value = 5313.5557
$ws.Range("K122").Value = 10988.7
$ws.Range("L122").Value = 15940.6671
$ws.Range("M122").Value = -8538.700000000001
$ws.Range("N122").Value = -20840.6671
$ws.Range("H126").Value = 58392.75
$ws.Range("J126").Value = 9761.75
$ws.Range("L126").Value = 29285.25
$ws.Range("N126").Value = -34225.25
$ws.Range("H132").Value = 8951.143
$ws.Range("I132").Value = 9797.166999999999
$ws.Range("J132").Value = 3875
$ws.Range("K132").Value = 29391.501
$ws.Range("L132").Value = 11625
$ws.Range("M132").Value = -26861.501
$ws.Range("N132").Value = -16685

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1361.1515
$ws.Range("I107").Value = 1233.6818
$ws.Range("K107").Value = 3701.0454
$ws.Range("M107").Value = -1781.0454
$ws.Range("H126").Value = 24478.523
$ws.Range("I126").Value = 48551.555
$ws.Range("K126").Value = 145654.665
$ws.Range("M126").Value = -143184.665
